$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 119; this pushes the existing rows
# 119-136 down to 121-138, matching the target layout exactly.
$ws.Rows("119:120").Insert()

# Populate the two newly inserted rows with the new weekly observation.
# Row 119 - "Primera"
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119, 3).Value = "Ñuble"
$ws.Cells.Item(119, 4).Value = 44918
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112040
$ws.Cells.Item(119, 7).Value = "Cilantro"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 400
$ws.Cells.Item(119, 11).Value = 600
$ws.Cells.Item(119, 12).Value = 700
$ws.Cells.Item(119, 13).Value = 650
$ws.Cells.Item(119, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(119, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(119, 16).Value = 650
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# Row 120 - "Segunda"
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44918
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112040
$ws.Cells.Item(120, 7).Value = "Cilantro"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Segunda"
$ws.Cells.Item(120, 10).Value = 300
$ws.Cells.Item(120, 11).Value = 500
$ws.Cells.Item(120, 12).Value = 500
$ws.Cells.Item(120, 13).Value = 500
$ws.Cells.Item(120, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(120, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(120, 16).Value = 500
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
